# Update the AccessKey test data value and adjust the active selection
# on the "APIData" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("APIData")
$ws.Activate()

# D2 (AccessKey column) changes from "c" to the new access key value.
$ws.Range("D2").Value = "7fe67bf08c80ded756e598d6f8fedaea"

# Move the active selection from B2 to A2.
$ws.Range("A2").Select()
